$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 6
$ws.Range("H2").Value = 7
$ws.Range("J2").Value = 3

# Row 3
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 10
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = 2

# Row 4
$ws.Range("B4").Value = 4
$ws.Range("D4").Value = 1
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 6

# Row 5
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 8
